$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format from an existing header cell (AC1) to the new
# header cells so they match the bold/bordered/centered header look.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add new header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record data for each player row (rows 2-41)
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 93   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 69   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
